# Policy.xlsx edit: append a trailing semicolon to the two action-cell
# DRL snippets in the "policy" RuleTable (row 8, columns C and D), and move
# the active selection to D9 (matches the author's final cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("policy")

$ws.Range("C8").Value = '$policy.setName("$param");'
$ws.Range("D8").Value = '$policy.setAmount($param);'

$ws.Range("D9").Select()
